# Append new log rows (273-283) to the "check_availability" log sheet,
# matching the rows produced by the application run on 2024-10-02.
#
# Every new cell is written as plain text (matching the existing rows,
# which are all inline/shared strings with no number formatting) -
# so date-looking ("2024-10-02") and currency-looking ("$199.99") values
# are forced to stay text instead of being auto-coerced into date serials
# or numbers by the normal Value setter. We do that by setting the
# NumberFormat to Text ("@") before assigning values, then clearing the
# formatting again afterwards so the cells end up with no explicit style,
# exactly like the rest of the log rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 273; A = "2024-10-02 16:54:42"; B = "check_availability"; C = "https://example.com";         D = "Checked availability: Selected or default date is available for booking.";                 E = "2024-10-02"; F = "16:54:42" },
    @{ Row = 274; A = "2024-10-02 17:15:04"; B = "check_availability"; C = "https://example.com";         D = "Checked availability: Selected or default date current date is available for booking.";    E = "2024-10-02"; F = "17:15:04" },
    @{ Row = 275; A = "2024-10-02 17:15:04"; B = "check_availability"; C = "https://example.com";         D = "Failed to check availability: Failed to check availability";                               E = "2024-10-02"; F = "17:15:04" },
    @{ Row = 276; A = "2024-10-02 17:15:04"; B = "check_availability"; C = "https://example.com";         D = "Checked availability: No availability for the selected date.";                             E = "2024-10-02"; F = "17:15:04" },
    @{ Row = 277; A = "2024-10-02 17:15:05"; B = "check_availability"; C = "https://example.com/product"; D = '$199.99';                                                                                     E = "2024-10-02"; F = "17:15:05" },
    @{ Row = 278; A = "2024-10-02 17:15:05"; B = "check_availability"; C = "invalid_url";                 D = "Error fetching price: Invalid URL";                                                         E = "2024-10-02"; F = "17:15:05" },
    @{ Row = 279; A = "2024-10-02 17:15:05"; B = "check_availability"; C = "https://example.com";         D = "Checked availability: Selected or default date is available for booking.";                 E = "2024-10-02"; F = "17:15:05" },
    @{ Row = 280; A = "2024-10-02 17:15:06"; B = "check_availability"; C = "https://example.com";         D = "Failed to check availability: Failed to check availability";                               E = "2024-10-02"; F = "17:15:06" },
    @{ Row = 281; A = "2024-10-02 17:15:07"; B = "check_availability"; C = "https://example.com/product"; D = "100 USD";                                                                                     E = "2024-10-02"; F = "17:15:07" },
    @{ Row = 282; A = "2024-10-02 18:03:26"; B = "check_availability"; C = "https://example.com/product"; D = "100 USD";                                                                                     E = "2024-10-02"; F = "18:03:26" },
    @{ Row = 283; A = "2024-10-02 18:03:55"; B = "check_availability"; C = "https://example.com/product"; D = "100 USD";                                                                                     E = "2024-10-02"; F = "18:03:55" }
)

foreach ($r in $rows) {
    $rowRange = $ws.Range("A$($r.Row):F$($r.Row)")
    # Force text storage for the whole row first so date/number-looking
    # strings (columns E and D in a couple of rows) are not auto-converted.
    $rowRange.NumberFormat = "@"

    $ws.Range("A$($r.Row)").Value = $r.A
    $ws.Range("B$($r.Row)").Value = $r.B
    $ws.Range("C$($r.Row)").Value = $r.C
    $ws.Range("D$($r.Row)").Value = $r.D
    $ws.Range("E$($r.Row)").Value = $r.E
    $ws.Range("F$($r.Row)").Value = $r.F

    # Drop the temporary Text format again so the new cells end up
    # unstyled, matching the rest of the data rows in this sheet.
    $rowRange.ClearFormats()
}
